$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PSSM values (rows 2-21, columns B-K), keyed by cell reference.
$values = @{
    "B2" = -15.00790327338205
    "C2" = 4.321886496609831
    "D2" = -15.00790327338205
    "E2" = -15.00790327338205
    "F2" = -15.00790327338205
    "G2" = -15.00790327338205
    "H2" = -15.00790327338205
    "I2" = -15.00790327338205
    "J2" = -15.00790327338205
    "K2" = -15.00790327338205
    "B3" = -15.00790327338205
    "C3" = -15.00790327338205
    "D3" = -15.00790327338205
    "E3" = -15.00790327338205
    "F3" = -15.00790327338205
    "G3" = -15.00790327338205
    "H3" = -15.00790327338205
    "I3" = 1.519222423966517
    "J3" = -15.00790327338205
    "K3" = -15.00790327338205
    "B4" = -15.00790327338205
    "C4" = -15.00790327338205
    "D4" = 1.986682195949108
    "E4" = -15.00790327338205
    "F4" = 3.575829600796953
    "G4" = -15.00790327338205
    "H4" = 1.83072208044706
    "I4" = -15.00790327338205
    "J4" = 2.503702694815377
    "K4" = -15.00790327338205
    "B5" = -15.00790327338205
    "C5" = -15.00790327338205
    "D5" = -15.00790327338205
    "E5" = -15.00790327338205
    "F5" = -15.00790327338205
    "G5" = 3.280354522118533
    "H5" = -15.00790327338205
    "I5" = -15.00790327338205
    "J5" = -15.00790327338205
    "K5" = -15.00790327338205
    "B6" = -15.00790327338205
    "C6" = -15.00790327338205
    "D6" = -15.00790327338205
    "E6" = -15.00790327338205
    "F6" = -15.00790327338205
    "G6" = -15.00790327338205
    "H6" = -15.00790327338205
    "I6" = -15.00790327338205
    "J6" = -15.00790327338205
    "K6" = -15.00790327338205
    "B7" = 3.014795047148663
    "C7" = -15.00790327338205
    "D7" = -15.00790327338205
    "E7" = -15.00790327338205
    "F7" = -15.00790327338205
    "G7" = -15.00790327338205
    "H7" = -15.00790327338205
    "I7" = -15.00790327338205
    "J7" = -15.00790327338205
    "K7" = -15.00790327338205
    "B8" = -15.00790327338205
    "C8" = -15.00790327338205
    "D8" = -15.00790327338205
    "E8" = 2.156974199016473
    "F8" = -15.00790327338205
    "G8" = -15.00790327338205
    "H8" = -15.00790327338205
    "I8" = -15.00790327338205
    "J8" = -15.00790327338205
    "K8" = -15.00790327338205
    "B9" = 3.574948029509103
    "C9" = -15.00790327338205
    "D9" = -15.00790327338205
    "E9" = -15.00790327338205
    "F9" = -15.00790327338205
    "G9" = -15.00790327338205
    "H9" = -15.00790327338205
    "I9" = -15.00790327338205
    "J9" = -15.00790327338205
    "K9" = -15.00790327338205
    "B10" = -15.00790327338205
    "C10" = -15.00790327338205
    "D10" = -15.00790327338205
    "E10" = -15.00790327338205
    "F10" = -15.00790327338205
    "G10" = -15.00790327338205
    "H10" = -15.00790327338205
    "I10" = 0.9877976405883303
    "J10" = -15.00790327338205
    "K10" = 2.038762479633886
    "B11" = -15.00790327338205
    "C11" = -15.00790327338205
    "D11" = -15.00790327338205
    "E11" = 2.475322069443552
    "F11" = -15.00790327338205
    "G11" = 2.06704304018123
    "H11" = -15.00790327338205
    "I11" = -15.00790327338205
    "J11" = -15.00790327338205
    "K11" = 1.620955171143142
    "B12" = -15.00790327338205
    "C12" = -15.00790327338205
    "D12" = -15.00790327338205
    "E12" = -15.00790327338205
    "F12" = -15.00790327338205
    "G12" = -15.00790327338205
    "H12" = -15.00790327338205
    "I12" = -15.00790327338205
    "J12" = -15.00790327338205
    "K12" = -15.00790327338205
    "B13" = -15.00790327338205
    "C13" = -15.00790327338205
    "D13" = -15.00790327338205
    "E13" = 2.091701375028513
    "F13" = -15.00790327338205
    "G13" = -15.00790327338205
    "H13" = -15.00790327338205
    "I13" = -15.00790327338205
    "J13" = 1.566705035591324
    "K13" = 2.442671688684698
    "B14" = -15.00790327338205
    "C14" = -15.00790327338205
    "D14" = 2.143162829544496
    "E14" = -15.00790327338205
    "F14" = -15.00790327338205
    "G14" = -15.00790327338205
    "H14" = -15.00790327338205
    "I14" = -15.00790327338205
    "J14" = -15.00790327338205
    "K14" = 1.89028476311486
    "B15" = -15.00790327338205
    "C15" = -15.00790327338205
    "D15" = 1.071018889690699
    "E15" = -15.00790327338205
    "F15" = -15.00790327338205
    "G15" = -15.00790327338205
    "H15" = -15.00790327338205
    "I15" = -15.00790327338205
    "J15" = -15.00790327338205
    "K15" = -15.00790327338205
    "B16" = -15.00790327338205
    "C16" = -15.00790327338205
    "D16" = -15.00790327338205
    "E16" = -15.00790327338205
    "F16" = -15.00790327338205
    "G16" = -15.00790327338205
    "H16" = -15.00790327338205
    "I16" = -15.00790327338205
    "J16" = 2.065997692271274
    "K16" = -15.00790327338205
    "B17" = -15.00790327338205
    "C17" = -15.00790327338205
    "D17" = 1.47485742447832
    "E17" = -15.00790327338205
    "F17" = -15.00790327338205
    "G17" = -15.00790327338205
    "H17" = 1.981812745571853
    "I17" = 0.9310560145115988
    "J17" = 1.928976168013336
    "K17" = -15.00790327338205
    "B18" = -15.00790327338205
    "C18" = -15.00790327338205
    "D18" = -15.00790327338205
    "E18" = -15.00790327338205
    "F18" = -15.00790327338205
    "G18" = -15.00790327338205
    "H18" = 2.09087933240493
    "I18" = 0.5948245844121357
    "J18" = 1.753089492269438
    "K18" = -15.00790327338205
    "B19" = -15.00790327338205
    "C19" = -15.00790327338205
    "D19" = 1.993295979007913
    "E19" = -15.00790327338205
    "F19" = -15.00790327338205
    "G19" = -15.00790327338205
    "H19" = 1.859001324985697
    "I19" = 1.523798005946399
    "J19" = -15.00790327338205
    "K19" = -15.00790327338205
    "B20" = -15.00790327338205
    "C20" = -15.00790327338205
    "D20" = 1.463152693840342
    "E20" = -15.00790327338205
    "F20" = 3.013494252451557
    "G20" = -15.00790327338205
    "H20" = 1.219054493356318
    "I20" = 3.146919779106767
    "J20" = -15.00790327338205
    "K20" = 1.876279172719423
    "B21" = -15.00790327338205
    "C21" = -15.00790327338205
    "D21" = -15.00790327338205
    "E21" = 2.51508429584475
    "F21" = -15.00790327338205
    "G21" = 2.607216238414694
    "H21" = 1.18710359984428
    "I21" = -15.00790327338205
    "J21" = -15.00790327338205
    "K21" = -15.00790327338205
}

foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
